$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "Humacao"
$ws.Range("B15").Value = 6080

$ws.Range("E5").Select()
